# Update export file project
# Rebuilds "Capstone 1" and "Capstone 2" sheets: adds a "Mentor" column,
# adds new data rows, and applies the report-style formatting (Times New
# Roman fonts, thin borders, centred/wrapped header, merged sub-columns on
# sheet1, column widths) that the refreshed export uses.

$wb = $excel.ActiveWorkbook

$headers = @("No.", "Student Code", "First name", "Last name", "Class", "Group", "Mentor", "Topic", "Description")

# Column widths (in the units Range.ColumnWidth expects) that render as
# 5 / 15 / 20 / 15 / 15 / 10 / 30 / 50 / 75 characters once saved - Excel
# adds ~0.8333 of padding to whatever ColumnWidth is set to.
$colWidths = @(4.166666666666667, 14.166666666666666, 19.166666666666668, 14.166666666666666, 14.166666666666666, 9.166666666666666, 29.166666666666668, 49.166666666666664, 74.16666666666667)

function Format-ReportSheet($ws) {
    # Column widths first (before any cell styling so no stray styled
    # blank cells get materialized).
    for ($i = 1; $i -le $colWidths.Length; $i++) {
        $ws.Columns.Item($i).ColumnWidth = $colWidths[$i - 1]
    }
}

function Write-HeaderRow($ws, $row) {
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $headers[$i]
    }
}

function Set-TextValue($cell, [string]$value) {
    # Forces a value that *looks* numeric (student codes, "123", ...) to be
    # stored as text instead of being auto-coerced to a number, without
    # leaving a lasting "@" number-format footprint on the cell's style.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Apply-HeaderStyle($rng) {
    $rng.Font.Name = "Times new roman"
    $rng.Font.Bold = $true
    $rng.Font.Size = 13
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

function Apply-DataStyle($rng) {
    $rng.Font.Name = "Times New Roman"
    $rng.Font.Size = 13
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
    $rng.HorizontalAlignment = -4131
    $rng.VerticalAlignment = -4108
    $rng.WrapText = $true
}

# ---------------------------------------------------------------------
# Sheet 1: "Capstone 1"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Capstone 1")
Format-ReportSheet $ws1

# Header row
Write-HeaderRow $ws1 1
$ws1.Rows.Item(1).RowHeight = 18
Apply-HeaderStyle $ws1.Range("A1:I1")

# Row 2
$ws1.Range("A2").Value = 1
$ws1.Range("C2").Value = "Nguyễn Văn"
$ws1.Range("D2").Value = "Tiên"
$ws1.Range("E2").Value = "CMU-TPM"
$ws1.Range("F2").Value = "C1SE.21"
$ws1.Range("G2").Value = "null null,null null"
$ws1.Range("H2").Value = "Capstone 123"
Set-TextValue $ws1.Range("B2") "24211202633"
Set-TextValue $ws1.Range("I2") "123"
$ws1.Rows.Item(2).RowHeight = 18
Apply-DataStyle $ws1.Range("A2:I2")

# Row 3 (F3:I3 are merged down into row 4)
$ws1.Range("A3").Value = 2
$ws1.Range("C3").Value = "Hà Đức"
$ws1.Range("D3").Value = "Phước"
$ws1.Range("E3").Value = "CMU-TPM"
$ws1.Range("F3").Value = "C1SE.23"
$ws1.Range("G3").Value = "null null"
$ws1.Range("H3").Value = "Capstone 123"
Set-TextValue $ws1.Range("B3") "24211202634"
Set-TextValue $ws1.Range("I3") "123"
$ws1.Rows.Item(3).RowHeight = 18

# Row 4
$ws1.Range("A4").Value = 3
$ws1.Range("C4").Value = "Nguyễn Văn"
$ws1.Range("D4").Value = "Tiên"
$ws1.Range("E4").Value = "CMU-TPM"
Set-TextValue $ws1.Range("B4") "24211208536"
$ws1.Rows.Item(4).RowHeight = 18

Apply-DataStyle $ws1.Range("A3:I4")

# Merge the Group / Topic / Description columns across rows 3-4
$ws1.Range("F3:F4").Merge()
$ws1.Range("G3:G4").Merge()
$ws1.Range("H3:H4").Merge()
$ws1.Range("I3:I4").Merge()

# ---------------------------------------------------------------------
# Sheet 2: "Capstone 2"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Capstone 2")
Format-ReportSheet $ws2

Write-HeaderRow $ws2 1
$ws2.Rows.Item(1).RowHeight = 18
Apply-HeaderStyle $ws2.Range("A1:I1")

$ws2.Range("A2").Value = 1
$ws2.Range("C2").Value = "Đậu Minh"
$ws2.Range("D2").Value = "Hoàng"
$ws2.Range("E2").Value = "CMU-TPM"
$ws2.Range("F2").Value = "C2SE.01"
$ws2.Range("G2").Value = ""
$ws2.Range("H2").Value = "Doctor"
$ws2.Range("I2").Value = "No"
Set-TextValue $ws2.Range("B2") "24211208533"
$ws2.Rows.Item(2).RowHeight = 18
Apply-DataStyle $ws2.Range("A2:I2")
